$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Test Case Name"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "testUrl"
$ws.Range("E1").Value = "stars_1"
$ws.Range("F1").Value = "stars_2"
$ws.Range("G1").Value = "policy"
$ws.Range("H1").Value = "reviewText"

# Update data row (row 2)
$ws.Range("A2").Value = "Add"
$ws.Range("B2").Value = "facebook_test@inbox.ru"
$ws.Range("C2").Value = "Wallet@123"
$ws.Range("D2").Value = "https://wallethub.com/profile/test_insurance_company/"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4"
$ws.Range("F2").Value = "5"
$ws.Range("G2").Value = "health"
$ws.Range("H2").Value = "Lorem ipsum dolor sit amet, consectetuer adipiscing elit. Aenean commodo ligula eget dolor. Aenean massa. Cum sociis natoque penatibus et magnis dis parturient montes, nascetur ridiculus mus. Donec qu"

# Remove hyperlink on B2 (previously existed)
$ws.Hyperlinks.Delete()

# Column widths
$ws.Columns.Item(2).ColumnWidth = 26.88671875
$ws.Columns.Item(3).ColumnWidth = 15.5546875
$ws.Columns.Item(4).ColumnWidth = 27.109375

# Select B2
$ws.Range("B2").Select()
